# Swap the species-observation data between row 2 and row 4 on the active sheet.
# Columns involved: A, B, E, F, G, H, Q, R, Z, AB
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R", "Z", "AB")

foreach ($col in $cols) {
    $rangeRow2 = $ws.Range($col + "2")
    $rangeRow4 = $ws.Range($col + "4")

    $val2 = $rangeRow2.Value()
    $val4 = $rangeRow4.Value()

    $rangeRow2.Value = $val4
    $rangeRow4.Value = $val2
}
